# MITHEx_inputs.xlsx - "Remove cycle efficiency calculation for air brayton
# cycle. Remove excessive print statements."
#
# The XLSX-visible part of this commit is just UI/state cleanup that was
# bundled into the same save as the (non-spreadsheet) code changes:
#   - Clear the stray "Secondary Mass Flow Rate (kg/s)" value (B10) that had
#     been left on the "Plant Description" sheet.
#   - Leave the workbook with "Plant Description" as the active/selected
#     sheet (instead of "Cycle Parameters"), with the various sheets'
#     last-used selections updated to where the author had clicked.

$wb = $excel.ActiveWorkbook

# --- HX Parameters: just a cursor move (A5 -> A8), no data changed ---
$wsHX = $wb.Worksheets.Item("HX Parameters")
$wsHX.Range("A8").Select()

# --- Cycle Parameters: no longer the active tab; selection (B3) unchanged ---
$wsCycle = $wb.Worksheets.Item("Cycle Parameters")
$wsCycle.Range("B3").Select()

# --- Plant Description: clear the leftover Secondary Mass Flow Rate value,
#     then make this the active sheet/selection (matches the saved file) ---
$wsPlant = $wb.Worksheets.Item("Plant Description")
$wsPlant.Range("B10").ClearContents()

$wsPlant.Activate()
$wsPlant.Range("B7").Select()

# Best-effort: also nudge the saved window position to the right (the author's
# workbook window moved from x=760 to x=12780 between saves).
$excel.ActiveWindow.Left = 12780
